# Add "shock_raw" and "extreme_level_raw" columns (E, F) holding the
# unrounded numeric values behind the formatted "shock" / "extreme_level"
# text columns (C, D), so the raw numbers survive for downstream use.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: copy the header formatting from C1/D1 onto E1/F1 ---
$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "shock_raw"

$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F1").Value = "extreme_level_raw"

$excel.CutCopyMode = 0

# --- Raw numeric values for column E (shock_raw) ---
$shockRaw = @{
    2  = 0.001283697047496757
    3  = 0.001254705144291046
    4  = 0.001226993865030668
    5  = 0.001200480192076858
    6  = 0.001175088131609936
    7  = 0.00115074798619097
    8  = 0.001127395715896373
    9  = 0.001104972375690627
    10 = 0.001083423618634782
    11 = 0.01381509032943673
    12 = 13
    13 = 13
    14 = 13
    15 = 13
    16 = 0.003535126912538544
    17 = 13
    18 = 13
    19 = 13
    20 = 0.0008382229673093988
    21 = 13
}

foreach ($r in $shockRaw.Keys) {
    $ws.Cells.Item($r, 5).Value = $shockRaw[$r]
}

# --- Raw numeric values for column F (extreme_level_raw); rows with no
#     numeric extreme level keep an empty (but present) cell ---
$extremeRaw = @{
    12 = 972
    13 = 990
    14 = 1008
    15 = 1026
    16 = 0.003573020113968983
    17 = 1152
    18 = 1170
    19 = 1188
    21 = 1224
}

foreach ($r in $extremeRaw.Keys) {
    $ws.Cells.Item($r, 6).Value = $extremeRaw[$r]
}

$emptyExtremeRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 20)
foreach ($r in $emptyExtremeRows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = "'"
    $cell.ClearFormats()
}
